# Add a "Save" column (H) to the s_vals sheet, matching the existing
# header style used by the other header cells (e.g. G1), and fill in the
# corresponding data value for row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1, bold/border/
# centered style) onto the new header cell H1 so it reuses the same style
# definition instead of creating a new one.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Set the new header text and the new data value.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0

Write-Output "Added Save column (H1:H2)"
